# Applies the "cryptos" price/volume refresh for Thu Aug  3 22:46:27 UTC 2023.
# Each cell is written to the exact final text shown in the source diff.
# Numeric-looking Price (column D) strings are written with a leading
# apostrophe so Excel stores them as literal text (preserving formats like
# "1.000" or "0.1390") instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.185.63'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '1.834.37'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("D5").Value = '''240.53'
$ws.Range("E5").Value = '  -0.32%  '

$ws.Range("D6").Value = '''0.6647'
$ws.Range("E6").Value = '  -3.36%  '

$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("D9").Value = '''0.2914'
$ws.Range("E9").Value = '  -2.83%  '

$ws.Range("D10").Value = '''22.61'
$ws.Range("E10").Value = '  -2.67%  '

$ws.Range("D11").Value = '''0.07692'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '1.833.88'
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").Value = '''4.968'
$ws.Range("E13").Value = '  -1.88%  '

$ws.Range("D14").Value = '''0.6643'
$ws.Range("E14").Value = '  -2.77%  '

$ws.Range("D15").Value = '''83.44'
$ws.Range("E15").Value = '  -4.42%  '

$ws.Range("D16").Value = '''6.078'
$ws.Range("E16").Value = '  -1.51%  '

$ws.Range("D17").Value = '29.191.47'
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").Value = '''0.000008253'
$ws.Range("E18").Value = '  +0.83%  '

$ws.Range("D19").Value = '''225.49'
$ws.Range("E19").Value = '  -1.63%  '

$ws.Range("D20").Value = '''12.42'
$ws.Range("E20").Value = '  -1.07%  '

$ws.Range("E21").Value = '  +0.13%  '

$ws.Range("D22").Value = '''7.121'
$ws.Range("E22").Value = '  -3.80%  '

$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = '''160.59'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").Value = '''8.623'
$ws.Range("E25").Value = '  -1.68%  '

$ws.Range("D26").Value = '''0.1390'
$ws.Range("E26").Value = '  -4.12%  '

$ws.Range("D27").Value = '''17.88'
$ws.Range("E27").Value = '  -1.21%  '

$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").Value = '''4.104'
$ws.Range("E29").Value = '  -4.16%  '

$ws.Range("D30").Value = '''4.023'
$ws.Range("E30").Value = '  -3.04%  '

$ws.Range("D31").Value = '''1.184'
$ws.Range("E31").Value = '  -1.17%  '

$ws.Range("D32").Value = '''0.05285'
$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").Value = '''1.864'
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("D34").Value = '''0.7486'
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("D35").Value = '''1.129'
$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("D36").Value = '''2.678'
$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("D37").Value = '1.301.80'
$ws.Range("E37").Value = '  -0.35%  '

$ws.Range("E38").Value = '  -2.14%  '

$ws.Range("D39").Value = '''2.717'
$ws.Range("E39").Value = '  -0.26%  '

$ws.Range("D40").Value = '''0.9196'
$ws.Range("E40").Value = '  -1.62%  '

$ws.Range("D41").Value = '''5.934'
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").Value = '''0.08505'
$ws.Range("E42").Value = '  +14.90%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.001'
$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''102.22'
$ws.Range("E44").Value = '  -2.83%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '''0.00000000128'
$ws.Range("E45").Value = '  +4.53%  '

$ws.Range("D46").Value = '1.981.45'
$ws.Range("E46").Value = '  -0.26%  '

$ws.Range("D47").Value = '''0.5161'
$ws.Range("E47").Value = '  -0.64%  '

$ws.Range("D48").Value = '''1.764'
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("D49").Value = '''63.20'
$ws.Range("E49").Value = '  -2.58%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05932'
$ws.Range("E50").Value = '  -0.33%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''8.989'
$ws.Range("E51").Value = '  -5.36%  '
